# Update "latest output (run 5)" optimisation results.
#
# Sheet "Schedule": row 2 gets new computed values, and the old row 3
# (second pumping block) is removed entirely.
#
# Sheet "Detailed": a new half-hourly record is inserted as row 2
# (shifting every following row down by one) and the Price / Type /
# Pump_Status values are refreshed to the latest optimisation run.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet 1: "Schedule"
# ---------------------------------------------------------------
$schedule = $wb.Worksheets.Item("Schedule")

# New values for row 2
$schedule.Cells.Item(2, 1).Value = 46037
$schedule.Cells.Item(2, 2).Value = 46037.66666666666
$schedule.Cells.Item(2, 3).Value = 16
$schedule.Cells.Item(2, 4).Value = 60.48
$schedule.Cells.Item(2, 5).Value = 1632.584050499999
$schedule.Cells.Item(2, 6).Value = 26.99378390376983

# Old row 3 no longer exists in the latest run - remove it entirely
# (this also shrinks the sheet dimension from A1:F3 to A1:F2).
$schedule.Rows.Item(3).Delete()

# ---------------------------------------------------------------
# Sheet 2: "Detailed"
# ---------------------------------------------------------------
$detailed = $wb.Worksheets.Item("Detailed")

# Final table (rows 2-49) for the "Detailed" sheet after the refresh.
# Columns: DateTime, Price, Type, Date, Pump_Status
$rows = @(
    @(2,  46037,               78,       "historical", 46037, "ON"),
    @(3,  46037.02083333334,   76.66624, "historical", 46037, "ON"),
    @(4,  46037.04166666666,   78,       "historical", 46037, "ON"),
    @(5,  46037.0625,          82.35556, "historical", 46037, "ON"),
    @(6,  46037.08333333334,   78,       "historical", 46037, "ON"),
    @(7,  46037.10416666666,   76.64994, "forecast",   46037, "ON"),
    @(8,  46037.125,           73.59164, "forecast",   46037, "ON"),
    @(9,  46037.14583333334,   72.6703,  "forecast",   46037, "ON"),
    @(10, 46037.16666666666,   64.37671, "forecast",   46037, "ON"),
    @(11, 46037.1875,          61.34927, "forecast",   46037, "ON"),
    @(12, 46037.20833333334,   63.83266, "forecast",   46037, "ON"),
    @(13, 46037.22916666666,   84.7901,  "forecast",   46037, "ON"),
    @(14, 46037.25,            80.06771000000001, "forecast", 46037, "ON"),
    @(15, 46037.27083333334,   78.25467999999999, "forecast", 46037, "ON"),
    @(16, 46037.29166666666,   52.16133, "forecast",   46037, "ON"),
    @(17, 46037.3125,          51.07513, "forecast",   46037, "ON"),
    @(18, 46037.33333333334,   49.7656,  "forecast",   46037, "ON"),
    @(19, 46037.35416666666,   36.06,    "forecast",   46037, "ON"),
    @(20, 46037.375,           41.92585, "forecast",   46037, "ON"),
    @(21, 46037.39583333334,   10.52921, "forecast",   46037, "ON"),
    @(22, 46037.41666666666,   11.75013, "forecast",   46037, "ON"),
    @(23, 46037.4375,          36.06,    "forecast",   46037, "ON"),
    @(24, 46037.45833333334,   0.51,     "forecast",   46037, "ON"),
    @(25, 46037.47916666666,   2.83675,  "forecast",   46037, "ON"),
    @(26, 46037.5,             36.06,    "forecast",   46037, "ON"),
    @(27, 46037.52083333334,   36.0601,  "forecast",   46037, "ON"),
    @(28, 46037.54166666666,   0.51,     "forecast",   46037, "ON"),
    @(29, 46037.5625,          36.0601,  "forecast",   46037, "ON"),
    @(30, 46037.58333333334,   52.11471, "forecast",   46037, "ON"),
    @(31, 46037.60416666666,   56.98,    "forecast",   46037, "ON"),
    @(32, 46037.625,           58.40146, "forecast",   46037, "ON"),
    @(33, 46037.64583333334,   56.98,    "forecast",   46037, "ON"),
    @(34, 46037.66666666666,   47.42517, "forecast",   46037, "OFF"),
    @(35, 46037.6875,          53.07603, "forecast",   46037, "OFF"),
    @(36, 46037.70833333334,   57.6972,  "forecast",   46037, "OFF"),
    @(37, 46037.72916666666,   21.07294, "forecast",   46037, "OFF"),
    @(38, 46037.75,            61.49051, "forecast",   46037, "OFF"),
    @(39, 46037.77083333334,   70.66426, "forecast",   46037, "OFF"),
    @(40, 46037.79166666666,   101.77225, "forecast",  46037, "OFF"),
    @(41, 46037.8125,          132.06252, "forecast",  46037, "OFF"),
    @(42, 46037.83333333334,   158.99,   "forecast",   46037, "OFF"),
    @(43, 46037.85416666666,   120.01,   "forecast",   46037, "OFF"),
    @(44, 46037.875,           85.95,    "forecast",   46037, "OFF"),
    @(45, 46037.89583333334,   77.14255, "forecast",   46037, "OFF"),
    @(46, 46037.91666666666,   78.00005, "forecast",   46037, "OFF"),
    @(47, 46037.9375,          64.99985, "forecast",   46037, "OFF"),
    @(48, 46037.95833333334,   71.40000000000001, "forecast", 46037, "OFF"),
    @(49, 46037.97916666666,   72.27782999999999, "forecast", 46037, "OFF")
)

# Row 49 is brand new - give it the same number formats as row 48
# (datetime for column A, date for column D) before filling it in.
$detailed.Range("A49").NumberFormat = $detailed.Range("A48").NumberFormat
$detailed.Range("D49").NumberFormat = $detailed.Range("D48").NumberFormat

foreach ($row in $rows) {
    $r = $row[0]
    $detailed.Cells.Item($r, 1).Value = $row[1]
    $detailed.Cells.Item($r, 2).Value = $row[2]
    $detailed.Cells.Item($r, 3).Value = $row[3]
    $detailed.Cells.Item($r, 4).Value = $row[4]
    $detailed.Cells.Item($r, 5).Value = $row[5]
}

Write-Host "Workbook updated."
